# Edit: change the title on the "Quality of Life - Education" slide so it
# reads "Quality of Life - Crime" instead, splitting the title text into two
# runs ("Quality of Life " and "- Crime") the way PowerPoint does when you
# select the trailing portion of a run and retype it.

$p = $ppt.ActivePresentation

$targetOld = "Quality of Life - Education"
$prefix    = "Quality of Life "
$newTail   = "- Crime"

$found = $false

for ($si = 1; $si -le $p.Slides.Count -and -not $found; $si++) {
    $slide = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $slide.Shapes.Count -and -not $found; $shi++) {
        $shape = $slide.Shapes.Item($shi)

        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $tr = $shape.TextFrame.TextRange

            if ($tr.Text -eq $targetOld) {
                # Grab just the "- Education" tail (everything after the
                # shared "Quality of Life " prefix) and retype it in place;
                # PowerPoint keeps the untouched prefix as its own run and
                # creates a new run for the replaced tail.
                $tailLen = $tr.Length - $prefix.Length
                $tail = $tr.Characters($prefix.Length + 1, $tailLen)
                $tail.Text = $newTail

                $found = $true
            }
        }
    }
}

if (-not $found) {
    throw "Could not find a shape with text '$targetOld'"
}
